$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.961.73'
$ws.Range("E2").Value = '  -0.88%  '

$ws.Range("D3").Value = '1.918.15'

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5054'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4043'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08338'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.103'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.63%  '

$ws.Range("D13").Value = '1.915.03'
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.402'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.52%  '

$ws.Range("E15").Value = '  -1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.14%  '

$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06512'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.16%  '

$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.948'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").Value = '29.997.15'
$ws.Range("E23").Value = '  -0.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.197'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.63%  '

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.137.67'
$ws.Range("E26").Value = '  +1.13%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.316'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.132'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.67%  '

$ws.Range("E32").Value = '  -1.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.958'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.828'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02452'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.404'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06418'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2150'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.06%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.196'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.75%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.723'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6469'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.212'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.222'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6038'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.638'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.209'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '78.95'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.131'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.79%  '
